$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.110.21"
$ws.Range("E2").Value = "  -0.16%  "

$ws.Range("D3").Value = "2.319.81"
$ws.Range("E3").Value = "  +0.54%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'303.43"
$ws.Range("E5").Value = "  +0.44%  "

$ws.Range("E6").Value = "  -0.20%  "

$ws.Range("E7").Value = "  +0.57%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  +2.20%  "

$ws.Range("D10").Value = "'36.11"
$ws.Range("E10").Value = "  +4.78%  "

$ws.Range("E11").Value = "  -0.85%  "

$ws.Range("E12").Value = "  -1.48%  "

$ws.Range("D13").Value = "'17.74"
$ws.Range("E13").Value = "  -1.77%  "

$ws.Range("D14").Value = "'6.92"
$ws.Range("E14").Value = "  +1.43%  "

$ws.Range("D15").Value = "2.680.91"
$ws.Range("E15").Value = "  +0.56%  "

$ws.Range("D16").Value = "2.317.35"
$ws.Range("E16").Value = "  -0.15%  "

$ws.Range("E17").Value = "  -2.64%  "

$ws.Range("D18").Value = "43.029.04"

$ws.Range("D19").Value = "'13.18"
$ws.Range("E19").Value = "  +4.74%  "

$ws.Range("E20").Value = "  +1.32%  "

$ws.Range("E21").Value = "  +0.19%  "

$ws.Range("D22").Value = "'68.27"
$ws.Range("E22").Value = "  +0.67%  "

$ws.Range("D23").Value = "'239.94"
$ws.Range("E23").Value = "  +1.19%  "

$ws.Range("E24").Value = "  -2.47%  "

$ws.Range("E25").Value = "  -0.31%  "

$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.65%  "

$ws.Range("D27").Value = "'25.49"
$ws.Range("E27").Value = "  +2.62%  "

$ws.Range("D28").Value = "'169.24"
$ws.Range("E28").Value = "  +0.44%  "

$ws.Range("D29").Value = "'34.12"
$ws.Range("E29").Value = "  -0.28%  "

$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").Value = "'9.19"
$ws.Range("E30").Value = "  +0.36%  "

$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").Value = "'2.05"
$ws.Range("E31").Value = "  -6.03%  "

$ws.Range("D32").Value = "'4.98"
$ws.Range("E32").Value = "  +9.17%  "

$ws.Range("D33").Value = "'5.18"
$ws.Range("E33").Value = "  +2.75%  "

$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "  -0.05%  "

$ws.Range("D35").Value = "'17.95"
$ws.Range("E35").Value = "  +4.95%  "

$ws.Range("E36").Value = "  -1.33%  "

$ws.Range("D37").Value = "'0.0698"
$ws.Range("E37").Value = "  +0.90%  "

$ws.Range("E38").Value = "  +1.37%  "

$ws.Range("E39").Value = "  -0.09%  "

$ws.Range("D40").Value = "'2.78"
$ws.Range("E40").Value = "  -1.15%  "

$ws.Range("E41").Value = "  +0.16%  "

$ws.Range("D42").Value = "1.993.10"
$ws.Range("E42").Value = "  -0.41%  "

$ws.Range("E43").Value = "  +1.26%  "

$ws.Range("D44").Value = "'2.23"
$ws.Range("E44").Value = "  -5.89%  "

$ws.Range("D45").Value = "'10.21"
$ws.Range("E45").Value = "  +0.65%  "

$ws.Range("D46").Value = "'17.57"
$ws.Range("E46").Value = "  -0.63%  "

$ws.Range("E47").Value = "  -0.05%  "

$ws.Range("D48").Value = "'76.23"
$ws.Range("E48").Value = "  +8.39%  "

$ws.Range("E49").Value = "  -2.12%  "

$ws.Range("D50").Value = "2.547.05"
$ws.Range("E50").Value = "  +0.49%  "

$ws.Range("E51").Value = "  +0.57%  "
